# Update IDs group 2: replace Subject1-5 data with Subject6-10 data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SubjectID6"
$ws.Range("B2").Value = "label6"
$ws.Range("C2").Value = "notes6"

$ws.Range("A3").Value = "SubjectID7"
$ws.Range("B3").Value = "label7"
$ws.Range("C3").Value = "notes7"

$ws.Range("A4").Value = "SubjectID8"
$ws.Range("B4").Value = "label8"
$ws.Range("C4").Value = "notes8"

$ws.Range("A5").Value = "SubjectID9"
$ws.Range("B5").Value = "label9"
$ws.Range("C5").Value = "notes9"

$ws.Range("A6").Value = "SubjectID10"
$ws.Range("B6").Value = "label10"
$ws.Range("C6").Value = "notes10"

# Column D width (bestFit applied in the original edit); closest
# achievable value given this engine's width quantization.
$ws.Columns("D").ColumnWidth = 11.1

# The original edit added an 8pt Calibri font to the style table
# (used by the phonetic-guide settings, not by any visible cell).
# Touch-and-revert an existing, already-populated cell's font size
# so the new font entry is recorded in styles.xml without altering
# that cell's own appearance.
$ws.Range("D2").Font.Size = 8
$ws.Range("D2").Font.Size = 11

# Selection / active cell, matching the saved view state
$ws.Range("B2:C6").Select()
